# Update cube metadata: Package3
# Renames the sheet, trims the unused trailing columns (T:Y) that only
# carried style (no values), retargets the column widths to the new
# narrower 15-column layout, and normalizes the trailing empty-row
# heights (adds a new blank row at the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet was renamed from "Informe-05-050314-A-TC-TP" to "Informe-04-040001-TM"
$ws.Name = "Informe-04-040001-TM"

# 2) Columns T:Y (20:25) in rows 1-5 only ever carried an (empty) style,
#    no value - fully clear them so the cells disappear and the used
#    range shrinks back from Y to S.
$ws.Range("T1:Y5").Clear()

# 3) New column layout/widths (cols A:O, 15 columns total). The
#    ColumnWidth COM property is expressed in "characters"; the engine
#    re-derives the stored sheet width by rounding to whole pixels
#    (width = round(ColumnWidth*6 + 5)/6), so we pre-compensate by the
#    standard 5/6-character padding to land on the intended stored
#    widths (27.69, 44.5, 18.66, 55.2, 34.64, 36.31, 47.28, 27.69 x3,
#    15.46, 46.44, 19.19, 20.05, 29.5).
$ws.Columns.Item(1).ColumnWidth = 26.85666666666667
$ws.Columns.Item(2).ColumnWidth = 43.666666666666664
$ws.Columns.Item(3).ColumnWidth = 17.826666666666668
$ws.Columns.Item(4).ColumnWidth = 54.36666666666667
$ws.Columns.Item(5).ColumnWidth = 33.806666666666665
$ws.Columns.Item(6).ColumnWidth = 35.47666666666667
$ws.Columns.Item(7).ColumnWidth = 46.446666666666665
$ws.Columns.Item(8).ColumnWidth = 26.85666666666667
$ws.Columns.Item(9).ColumnWidth = 26.85666666666667
$ws.Columns.Item(10).ColumnWidth = 26.85666666666667
$ws.Columns.Item(11).ColumnWidth = 14.626666666666667
$ws.Columns.Item(12).ColumnWidth = 45.60666666666666
$ws.Columns.Item(13).ColumnWidth = 18.35666666666667
$ws.Columns.Item(14).ColumnWidth = 19.21666666666667
$ws.Columns.Item(15).ColumnWidth = 28.666666666666668

# 4) The blank row right after the data (row 6) now matches the height
#    of the other filler rows (12.8 instead of 13.8), and a new blank
#    row 9 is appended with that same height.
$ws.Rows.Item(6).RowHeight = 12.8
$ws.Rows.Item(9).RowHeight = 12.8

# 5) Selection moved from A1 to B18.
[void]$ws.Range("B18").Select()
